# Update period labels: swap the "2401" / "2311" text values
# so that the three rows read 2311, 2312, 2401 (previously 2401, 2312, 2311)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2311"
$ws.Range("E18").Value = "2401"

# Update "Salario Basico" (column G) amounts for the three rows
$ws.Range("G16").Value = 1300000
$ws.Range("G17").Value = 1300000
$ws.Range("G18").Value = 1300000
